# "Generate Report for Handoff"
#
# The localization status report is regenerated: the "In Translation"
# status becomes "Ready for handoff" everywhere it appears, and the two
# "Latest Handoff Datetime" timestamps associated with that status move
# forward a little (12:38:27 -> 12:39:11 and 12:38:22 -> 12:39:01).
#
# Because the new status text ("Ready for handoff") is longer than the
# old one ("In Translation"), Excel's column autosize bumped the Status
# columns a bit wider on every sheet that shows it.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-18 12:39:11"

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-18 12:39:01"

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-18 12:39:11"

# --- Widen the Status columns to fit the new, longer text --------------
# (mirrors Excel auto-fitting the column after the longer value landed)
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # column F
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33        # column C
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33        # column C
